# khl/Injuries_Master_Clubs.xlsx - refresh scrape run (2025-12-20 ~12:10-12:11 UTC)
#   - two previously-injured players recovered and drop off the "snapshot"
#     sheet; they're logged on the "returned" sheet instead
#   - every remaining "snapshot" row gets a fresh scraped_at timestamp

$wb = $excel.ActiveWorkbook

# --- Sheet "snapshot": remove recovered players, refresh scrape timestamps ---
$snap = $wb.Worksheets.Item("snapshot")

# Delete rows bottom-up so the row index of earlier rows stays valid
$snap.Rows.Item(25).Delete()   # СОЧ / Бикмуллин Рафаэль (recovered)
$snap.Rows.Item(14).Delete()   # НХК / Пастухов Илья (recovered)

# Refresh scraped_at (column K) for every remaining data row
$snap.Cells.Item(2, 11).Value = "2025-12-20T12:10:41.325078+00:00"
$snap.Cells.Item(3, 11).Value = "2025-12-20T12:10:43.876190+00:00"
$snap.Cells.Item(4, 11).Value = "2025-12-20T12:10:43.876208+00:00"
$snap.Cells.Item(5, 11).Value = "2025-12-20T12:10:43.876216+00:00"
$snap.Cells.Item(6, 11).Value = "2025-12-20T12:10:46.017382+00:00"
$snap.Cells.Item(7, 11).Value = "2025-12-20T12:10:48.247126+00:00"
$snap.Cells.Item(8, 11).Value = "2025-12-20T12:10:51.028642+00:00"
$snap.Cells.Item(9, 11).Value = "2025-12-20T12:10:51.028670+00:00"
$snap.Cells.Item(10, 11).Value = "2025-12-20T12:10:53.244768+00:00"
$snap.Cells.Item(11, 11).Value = "2025-12-20T12:10:58.685422+00:00"
$snap.Cells.Item(12, 11).Value = "2025-12-20T12:11:01.357328+00:00"
$snap.Cells.Item(13, 11).Value = "2025-12-20T12:11:03.596007+00:00"
$snap.Cells.Item(14, 11).Value = "2025-12-20T12:11:08.613617+00:00"
$snap.Cells.Item(15, 11).Value = "2025-12-20T12:11:08.613646+00:00"
$snap.Cells.Item(16, 11).Value = "2025-12-20T12:11:08.613664+00:00"
$snap.Cells.Item(17, 11).Value = "2025-12-20T12:11:08.613680+00:00"
$snap.Cells.Item(18, 11).Value = "2025-12-20T12:11:11.424003+00:00"
$snap.Cells.Item(19, 11).Value = "2025-12-20T12:11:11.424033+00:00"
$snap.Cells.Item(20, 11).Value = "2025-12-20T12:11:11.424049+00:00"
$snap.Cells.Item(21, 11).Value = "2025-12-20T12:11:14.419892+00:00"
$snap.Cells.Item(22, 11).Value = "2025-12-20T12:11:14.419925+00:00"
$snap.Cells.Item(23, 11).Value = "2025-12-20T12:11:14.419942+00:00"
$snap.Cells.Item(24, 11).Value = "2025-12-20T12:11:17.010146+00:00"
$snap.Cells.Item(25, 11).Value = "2025-12-20T12:11:17.010180+00:00"
$snap.Cells.Item(26, 11).Value = "2025-12-20T12:11:19.787048+00:00"
$snap.Cells.Item(27, 11).Value = "2025-12-20T12:11:19.787078+00:00"
$snap.Cells.Item(28, 11).Value = "2025-12-20T12:11:19.787096+00:00"
$snap.Cells.Item(29, 11).Value = "2025-12-20T12:11:21.988714+00:00"
$snap.Cells.Item(30, 11).Value = "2025-12-20T12:11:24.356544+00:00"
$snap.Cells.Item(31, 11).Value = "2025-12-20T12:11:24.356571+00:00"
$snap.Cells.Item(32, 11).Value = "2025-12-20T12:11:29.329726+00:00"
$snap.Cells.Item(33, 11).Value = "2025-12-20T12:11:29.329759+00:00"
$snap.Cells.Item(34, 11).Value = "2025-12-20T12:11:31.839144+00:00"
$snap.Cells.Item(35, 11).Value = "2025-12-20T12:11:31.839175+00:00"

# --- Sheet "returned": log the two recovered players ---
$ret = $wb.Worksheets.Item("returned")

$ret.Cells.Item(3, 1).Value = "НХК"
$ret.Cells.Item(3, 2).Value = "Нефтехимик"
$ret.Cells.Item(3, 3).Value = "Пастухов Илья"
$ret.Cells.Item(3, 4).Value = "1369_НХК_пастуховилья"
$ret.Cells.Item(3, 5).Value = "RETURN"
$ret.Cells.Item(3, 6).Value = "2025-12-20T20:11:32.343828+08:00"
# Column G ("changed_day") is plain text "2025-12-20" in the source data, not
# a real date - force text interpretation so Excel doesn't coerce it into a
# date serial, then drop the now-unneeded explicit number format again.
$ret.Cells.Item(3, 7).NumberFormat = "@"
$ret.Cells.Item(3, 7).Value = "2025-12-20"
$ret.Cells.Item(3, 7).ClearFormats()

$ret.Cells.Item(4, 1).Value = "СОЧ"
$ret.Cells.Item(4, 2).Value = "ХК Сочи"
$ret.Cells.Item(4, 3).Value = "Бикмуллин Рафаэль"
$ret.Cells.Item(4, 4).Value = "1369_СОЧ_бикмуллинрафаэль"
$ret.Cells.Item(4, 5).Value = "RETURN"
$ret.Cells.Item(4, 6).Value = "2025-12-20T20:11:32.343828+08:00"
$ret.Cells.Item(4, 7).NumberFormat = "@"
$ret.Cells.Item(4, 7).Value = "2025-12-20"
$ret.Cells.Item(4, 7).ClearFormats()
